$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.0343
$ws.Range("B3").Value = 0.0269
$ws.Range("C3").Value = -0.1016
$ws.Range("B4").Value = -0.2673
$ws.Range("C4").Value = -0.2714
$ws.Range("B5").Value = -0.4034
$ws.Range("C5").Value = -0.3736
$ws.Range("B6").Value = 0.0065
$ws.Range("C6").Value = -0.0015
$ws.Range("B7").Value = 0.3413
$ws.Range("C7").Value = 0.2624
$ws.Range("B8").Value = 0.1623
$ws.Range("C8").Value = 0.1604
$ws.Range("B9").Value = 0.1351
$ws.Range("C9").Value = 0.1508
$ws.Range("B10").Value = 0.8549
$ws.Range("C10").Value = 0.7927
$ws.Range("B11").Value = 0.4819
$ws.Range("C11").Value = 0.4929
$ws.Range("B12").Value = 0.5441
$ws.Range("C12").Value = 0.5229
$ws.Range("B13").Value = 0.3316
$ws.Range("C13").Value = 0.3228
$ws.Range("B14").Value = 0.0187
$ws.Range("C14").Value = -0.0193
$ws.Range("B15").Value = -0.2525
$ws.Range("C15").Value = -0.2237
$ws.Range("B16").Value = -0.2697
$ws.Range("C16").Value = -0.2749
$ws.Range("B17").Value = -0.4474
$ws.Range("C17").Value = -0.4232
$ws.Range("B18").Value = -0.2039
$ws.Range("C18").Value = -0.1912
$ws.Range("B19").Value = 0.4652
$ws.Range("C19").Value = 0.4139
$ws.Range("B20").Value = 0.09569999999999999
$ws.Range("C20").Value = 0.08
$ws.Range("B21").Value = -0.5489000000000001
$ws.Range("C21").Value = -0.5022
$ws.Range("B22").Value = -0.6326000000000001
$ws.Range("C22").Value = -0.6106
$ws.Range("B23").Value = -0.6268
$ws.Range("C23").Value = -0.5702
$ws.Range("B24").Value = 0.2522
$ws.Range("C24").Value = 0.2235
$ws.Range("B27").Value = 0.5477
$ws.Range("B29").Value = -0.0164
$ws.Range("B30").Value = -0.0663
$ws.Range("B31").Value = -0.402
$ws.Range("B32").Value = -0.9175
$ws.Range("B33").Value = -1.3759
$ws.Range("B34").Value = -1.5161
$ws.Range("B35").Value = -1.8329
$ws.Range("B36").Value = -1.1538
